# Apply the mappings.xlsx update:
#  - SupplierMappings: add a new raw-value row mapping "World Kinect Corporation"
#    to the same standardized value/domain as "World Kinect".
#  - TerminalMappings: add four new cleaned/capitalized terminal mapping rows.
#  - Re-select SupplierMappings as the active sheet/tab (it was previously
#    TerminalMappings), and update the remembered selections on both sheets.

$wb = $excel.ActiveWorkbook

# --- SupplierMappings: append a new row (row 8) ---
$wsSupplier = $wb.Worksheets.Item("SupplierMappings")
$wsSupplier.Cells.Item(8, 1).Value = "World Kinect Corporation"
$wsSupplier.Cells.Item(8, 2).Value = "World Fuels"
$wsSupplier.Cells.Item(8, 3).Value = "world-kinect.com"

# --- TerminalMappings: append four new rows (rows 52-55) ---
$wsTerminal = $wb.Worksheets.Item("TerminalMappings")
$newTerminalRows = @(
    @("SINCLAIR-HEP-KANSAS CITY-KS", "Kansas City KS Sinclair"),
    @("IL Wood River - KMEP", "Hartford IL Kinder Morgan"),
    @("IL, Cahokia, PSX", "Cahokia IL Phillips 66"),
    @("IA, Bettendorf, MG, Marathon", "Bettendorf IA Magellan")
)
$startRow = 52
for ($i = 0; $i -lt $newTerminalRows.Count; $i++) {
    $r = $startRow + $i
    $wsTerminal.Cells.Item($r, 1).Value = $newTerminalRows[$i][0]
    $wsTerminal.Cells.Item($r, 2).Value = $newTerminalRows[$i][1]
}

# --- Update selections / active sheet ---
$wsTerminal.Range("A11").Select() | Out-Null

$wsSupplier.Activate() | Out-Null
$wsSupplier.Range("B8").Select() | Out-Null

Write-Host "Applied mappings.xlsx update"
